# Logbook Update - 19th June
#
# The paragraph that used to hold only the "_GoBack" bookmark (and was
# bold via its paragraph-mark formatting) is expanded into four
# paragraphs:
#   1. A new diary entry about the split-button / font-awesome work.
#   2. A blank spacer paragraph.
#   3. A new bold "Monday 19th June 4pm " heading paragraph.
#   4. The original (now non-bold) paragraph that still carries the
#      "_GoBack" bookmark.

$d = $word.ActiveDocument

# Anchor: the existing "Monday 19th June 12pm" heading paragraph. The
# diff leaves this paragraph untouched and replaces the paragraph right
# after it (which used to hold only the "_GoBack" bookmark) with four
# new paragraphs, so we insert the new material right after this
# heading.
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "Monday 19th June 12pm") {
        $anchorPara = $cand
        break
    }
}

$anchorRange = $anchorPara.Range
$anchorRange.InsertParagraphAfter()

# --- Paragraph 1: the new diary entry -----------------------------------
$entryPara = $anchorPara.Next()
$pos = $entryPara.Range.Start

$chunk1 = "So I may have been using Foundation" + [char]0x2019 + "s split buttons to "
$chunk2 = "implement the call-to-action buttons on the first section of the site. "
$chunk3 = "I" + [char]0x2019 + "ve looked into using the font-awesome ones and changed the angle to " + [char]0x2018 + "0" + [char]0x2019 + "."
$chunk4 = " I needed to change some styling and HTML to implement it correctly. This took around 30-40 minutes. "

$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk1)
$pos = $pos + $chunk1.Length

$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk2)
$pos = $pos + $chunk2.Length

$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk3)
$pos = $pos + $chunk3.Length

$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk4)
$pos = $pos + $chunk4.Length

# --- Paragraph 2: blank spacer paragraph ---------------------------------
$entryRange = $entryPara.Range
$entryRange.InsertParagraphAfter()

# --- Paragraph 3: bold "Monday 19th June 4pm " heading -------------------
$blankPara = $entryPara.Next()
$blankRange = $blankPara.Range
$blankRange.InsertParagraphAfter()

$headingPara = $blankPara.Next()
$hPos = $headingPara.Range.Start

$h1 = "Monday 19"
$h2 = "th"
$h3 = " June 4pm "

$ip = $d.Range($hPos, $hPos)
$ip.InsertAfter($h1)
$r1 = $d.Range($hPos, $hPos + $h1.Length)
$r1.Font.Bold = $true

$p2start = $hPos + $h1.Length
$ip = $d.Range($p2start, $p2start)
$ip.InsertAfter($h2)
$r2 = $d.Range($p2start, $p2start + $h2.Length)
$r2.Font.Bold = $true
$r2.Font.Superscript = $true

$p3start = $p2start + $h2.Length
$ip = $d.Range($p3start, $p3start)
$ip.InsertAfter($h3)
$r3 = $d.Range($p3start, $p3start + $h3.Length)
$r3.Font.Bold = $true

# Bold the paragraph mark itself, so the heading paragraph's <w:pPr>
# carries <w:rPr><w:b/></w:rPr>, matching the other date headings.
$fullHeading = $headingPara.Range
$fullHeading.Font.Bold = $true

# --- Paragraph 4: the original bookmark paragraph loses its bold mark ---
$bookmarkPara = $headingPara.Next()
$bookmarkPara.Range.Font.Bold = $false

Write-Output "Logbook updated with the 19th June 4pm entry."
